$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary per-row (D, L, M, N, O, P, Q, R, S, T) get reshuffled across
# rows 2..18 according to the mapping below: newRow[col] = oldRow[mapping[newRow]][col]
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# Snapshot the current ("before") values for every row/column we touch.
$snapshot = @{}
foreach ($r in 2..18) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# Target row -> source row mapping (data permutation).
$mapping = @{
    2 = 5
    3 = 6
    4 = 2
    5 = 13
    6 = 10
    7 = 11
    8 = 17
    9 = 18
    10 = 8
    11 = 7
    12 = 15
    13 = 14
    14 = 16
    15 = 12
    16 = 9
    17 = 3
    18 = 4
}

foreach ($destRow in 2..18) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
